$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41. Existing rows 41..67 (and their data/style)
# shift down to 42..68, matching the diff's observed row-by-row shift where
# every subsequent row's content equals the prior row's pre-edit content.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly price record
# (date 2022-08-04, serial 44777) using the same constant columns shared by
# every other row in this sheet (A, B, C, E, F, G, H, I, N, O, Q, R).
$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Femacal de La Calera"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44777
$ws.Range("E41").Value = 5
$ws.Range("F41").Value = 100112035
$ws.Range("G41").Value = "Bruselas (repollito)"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 85
$ws.Range("K41").Value = 14500
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = 14735
$ws.Range("N41").Value = "`$/malla 15 kilos"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 982
$ws.Range("Q41").Value = 15
$ws.Range("R41").Value = "Hortaliza"
